# Correção nos dados e início da análise PNAD 2009
#
# The sheet had two "category header" rows with no data of their own
# ("situação do domicílio" in row 5 and "grandes regiões" in row 8),
# which were followed by the real data rows. The fix removes those two
# empty header rows entirely so that the remaining rows (urbana, rural,
# norte, nordeste, sudeste, sul) shift up and sit directly under
# "brasil", each keeping its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the empty "situação do domicílio" header row (row 5).
# This shifts urbana/rural/grandes regiões/norte/nordeste/sudeste/sul
# up by one row (old row 6 -> new row 5, etc.).
$ws.Rows("5").Delete()

# After the first deletion, the empty "grandes regiões" header row that
# used to be row 8 is now row 7. Remove it too, shifting norte/nordeste/
# sudeste/sul up by one more row.
$ws.Rows("7").Delete()
